$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - add new column I
$ws.Range("I1").Value = "Primary Hazard"

# Barcode column (E) must stay text, not be auto-converted to a number.
# Temporarily mark the range as Text so the numeric-looking barcode values
# are stored as shared strings instead of numbers, then restore the default
# "Normal" style so no residual cell formatting remains on these cells.
$ws.Range("E2:E4").NumberFormat = "@"

# Row 2 - update existing row to new hazard-report data
$ws.Range("A2").Value = "Science"
$ws.Range("B2").Value = "Basement"
$ws.Range("C2").Value = "13C"
$ws.Range("D2").Value = "13C"
$ws.Range("E2").Value = "17050002"
$ws.Range("F2").Value = "BASE CHEMICAL"
$ws.Range("G2").Value = 3223
$ws.Range("H2").Value = "gram (g)"
$ws.Range("I2").Value = "Base"

# Row 3 - new row
$ws.Range("A3").Value = "Science"
$ws.Range("B3").Value = "Basement"
$ws.Range("C3").Value = "13C"
$ws.Range("D3").Value = "13C"
$ws.Range("E3").Value = "17050007"
$ws.Range("F3").Value = "OXIDIZER CHEMICAL"
$ws.Range("G3").Value = 33
$ws.Range("H3").Value = "gram (g)"
$ws.Range("I3").Value = "Oxidizer"

# Row 4 - new row
$ws.Range("A4").Value = "Science"
$ws.Range("B4").Value = "Basement"
$ws.Range("C4").Value = "13C"
$ws.Range("D4").Value = "13C"
$ws.Range("E4").Value = "17050008"
$ws.Range("F4").Value = "REACTIVE CHEMICAL"
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = "gram (g)"
$ws.Range("I4").Value = "Reactive"

# Restore default styling on the barcode column now that the text values
# are locked in, so no extra formatting lingers on these cells.
$ws.Range("E2:E4").Style = "Normal"
